# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps produced by a fresh handback run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 6aa6b236-...md
$wsOverview.Range("G2").Value = "2016-08-29 23:08:45"

# zh-cn sheet (row 2 = 6aa6b236-...md entry)
$wsZhCn.Range("H2").Value = "2016-08-29 23:08:41"   # Correspond Handoff Datetime
$wsZhCn.Range("K2").Value = "2016-08-29 23:08:59"   # Correspond Handback DateTime

# de-de sheet (row 2 = 6aa6b236-...md entry)
$wsDeDe.Range("H2").Value = "2016-08-29 23:08:45"   # Correspond Handoff Datetime
$wsDeDe.Range("K2").Value = "2016-08-29 23:09:11"   # Correspond Handback DateTime
